$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Split the combined "Events" value currently in G2 into two separate events ---
$combined = $ws.Cells.Item(2, 7).Value()
$parts = $combined.Split(";")
$event1 = $parts[0]
$event2 = $parts[1]

# --- Prepare the new column H cells that actually hold data (H1 & H2),
#     copying only their own row's formatting from column G ---
$ws.Cells.Item(1, 7).Copy()
$ws.Cells.Item(1, 8).PasteSpecial(-4122)
$ws.Cells.Item(2, 7).Copy()
$ws.Cells.Item(2, 8).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Rename G1 header to "Event 1" and add new H1 header "Event 2" ---
# (write in the same order the values first appear so the shared-string
# table is populated in the same sequence as the authored workbook)
$ws.Cells.Item(1, 7).Value = "Event 1"
$ws.Cells.Item(2, 7).Value = $event1
$ws.Cells.Item(2, 8).Value = $event2
$ws.Cells.Item(1, 8).Value = "Event 2"

# --- Set the width of the new column H (closest achievable value to the
#     authored 31.44140625 OOXML column width, given COM's pixel rounding) ---
$ws.Columns.Item(8).ColumnWidth = 30.6

# --- Move the visible selection, mirroring the end-user's final click position ---
$ws.Application.Goto($ws.Range("F1"))
$ws.Range("G2").Select()
